$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") was refreshed for every data row (2-176):
# the "last changed" date serial moved from 45175 (2023-09-06) to
# 45177 (2023-09-08).
for ($r = 2; $r -le 176; $r++) {
    $ws.Cells.Item($r, 3).Value = 45177
}

# Row 5 (A 13746-2023) additionally lost one signal species ("Vedticka")
# from its species list, so the counts that summarize that list shrink
# by one each:
#   I5 (Signalarter)  5 -> 4
#   Q5 (Alla arter)  12 -> 11
# and the "Vedticka" line is removed from the R5 species-name list.
$ws.Range("I5").Value = 4
$ws.Range("Q5").Value = 11

$ws.Range("R5").Value = "Knärot`r`nKoralltaggsvamp`r`nOrange taggsvamp`r`nTalltita`r`nUllticka`r`nBrandticka`r`nBronshjon`r`nSvavelriska`r`nThomsons trägnagare`r`nKopparödla`r`nBlåsippa"
